# FAST_holdings.xlsx update: bump the "as of" date in the confidentiality
# notice (A13) from 2021-05-27 to 2021-05-28 and refresh the Weight /
# Percent Change figures for rows 2-10 (columns D and E) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; unprotect so the cells can be written, then
# restore protection afterwards.
$ws.Unprotect()

$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-28 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.09335771520324111
$ws.Range("E2").Value = 0.005164146071560349

$ws.Range("D3").Value = 0.1074359932747873
$ws.Range("E3").Value = 0

$ws.Range("D4").Value = 0.118667186240567
$ws.Range("E4").Value = 0.001164867593383523

$ws.Range("D5").Value = 0.1401386701497959
$ws.Range("E5").Value = 0.001979289244528282

$ws.Range("D6").Value = 0.1344613140694953
$ws.Range("E6").Value = 0.003595132743362761

$ws.Range("D7").Value = 0.1459987182329083
$ws.Range("E7").Value = -0.002371982698479158

$ws.Range("D8").Value = 0.1285116984027035
$ws.Range("E8").Value = 0.002893518518518379

$ws.Range("D9").Value = 0.1314287044265018
$ws.Range("E9").Value = -0.002599399145660208

$ws.Range("E10").Value = 0.001065034556887445

$ws.Protect()
